$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "Nicola Lorenzi"
$ws.Range("B4").Value = "Stefano Tita | Clitoriders"
$ws.Range("C4").Value = "Carlo Stedile | MAI UNA GIOIA"
$ws.Range("D4").Value = "Federico  Mortillaro | Clitoriders"
$ws.Range("E4").Value = "Leonardo  Parisi  | MediaserT"
$ws.Range("F4").Value = "Mattia Baldessarini | Shark Attack"
